# Auto-generated edit script applying cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.695.17"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.626.01"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.51"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.07"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.625.01"
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.417"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.243.75"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000209"
$ws.Range("E14").Value = "  +0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.98"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.601.72"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.807.36"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.61"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.15"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "428.47"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.621"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.90"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +4.11%  "
$ws.Range("E27").Value = "  +5.58%  "
$ws.Range("E28").Value = "  +4.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.52"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.625.29"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.159"
$ws.Range("E33").Value = "  +4.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.48"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.90"
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.69"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "176.93"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0863"
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.902"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.59"
$ws.Range("E45").Value = "  +8.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.13"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.18"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.94"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.961"
$ws.Range("E51").Value = "  +1.54%  "
